$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend styled region (fill/border) down to new rows 26-35 by copying format from row 25
$ws.Range("A25:F25").Copy()
$ws.Range("A26:F35").PasteSpecial(-4122)

# Populate data rows 3-35
# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Margarita"
$ws.Range("C3").Value = "Đồ Uống Nóng"
$ws.Range("D3").Value = 180000
$ws.Range("E3").Value = "Còn"
$ws.Range("F3").Value = "Margarita cổ điển với một chút chanh."
# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Chicken Wings"
$ws.Range("C4").Value = "Rượu Vang"
$ws.Range("D4").Value = 240003
$ws.Range("E4").Value = "Còn"
$ws.Range("F4").Value = "Cánh gà cay nồng."
# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Vegetarian Pizza"
$ws.Range("C5").Value = "Đồ Uống Lạnh"
$ws.Range("D5").Value = 312000
$ws.Range("E5").Value = "Còn"
$ws.Range("F5").Value = "Pizza chay ngon với rau sống tươi."
# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Seafood Pasta"
$ws.Range("C6").Value = "Đồ Uống Nóng"
$ws.Range("D6").Value = 384000
$ws.Range("E6").Value = "Còn"
$ws.Range("F6").Value = "Mì hải sản phong cách Ý với tỏi và thảo mộc.gg"
# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Chocolate Martini"
$ws.Range("C7").Value = "Nước Ép Trái Cây"
$ws.Range("D7").Value = 4444
$ws.Range("E7").Value = "Còn"
# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Beef Tacos"
$ws.Range("C8").Value = "Đồ Uống Nóng"
$ws.Range("D8").Value = 252003
$ws.Range("E8").Value = "Còn"
$ws.Range("F8").Value = "Tacos bò thơm ngon với salsa và guacamole."
# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Caprese Salad"
$ws.Range("C9").Value = "Đồ Uống Nóng"
$ws.Range("D9").Value = 168000
$ws.Range("E9").Value = "Còn"
$ws.Range("F9").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Caprese Salad"
$ws.Range("C10").Value = "Đồ Uống Nóng"
$ws.Range("D10").Value = 168000
$ws.Range("E10").Value = "Còn"
$ws.Range("F10").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Caprese Salad"
$ws.Range("C11").Value = "Đồ Uống Nóng"
$ws.Range("D11").Value = 168000
$ws.Range("E11").Value = "Còn"
$ws.Range("F11").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Caprese Salad"
$ws.Range("C12").Value = "Đồ Uống Nóng"
$ws.Range("D12").Value = 168000
$ws.Range("E12").Value = "Còn"
$ws.Range("F12").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Caprese Salad"
$ws.Range("C13").Value = "Đồ Uống Nóng"
$ws.Range("D13").Value = 168000
$ws.Range("E13").Value = "Còn"
$ws.Range("F13").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Caprese Salad"
$ws.Range("C14").Value = "Đồ Uống Nóng"
$ws.Range("D14").Value = 168000
$ws.Range("E14").Value = "Còn"
$ws.Range("F14").Value = "Salad Caprese cổ điển với cà chua và phô mai tươi."
# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Nước Chuối"
$ws.Range("C15").Value = "Nước Ép Trái Cây"
$ws.Range("D15").Value = 10000
$ws.Range("E15").Value = "Hết"
# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Bia Hà Nội"
$ws.Range("C16").Value = "Bia"
$ws.Range("D16").Value = 30000
$ws.Range("E16").Value = "Còn"
# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Mật ong"
$ws.Range("C17").Value = "Đồ Uống Nóng"
$ws.Range("D17").Value = 6000
$ws.Range("E17").Value = "Còn"
$ws.Range("F17").Value = "mật ong vị cà chua"
# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Shrimp Scampi"
$ws.Range("C18").Value = "Đồ Uống Nóng"
$ws.Range("D18").Value = 355203
$ws.Range("E18").Value = "Còn"
$ws.Range("F18").Value = "Shrimp scampi với bơ tỏi và gừng."
# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Beef Tacos"
$ws.Range("C19").Value = "Đồ Uống Nóng"
$ws.Range("D19").Value = 252001
$ws.Range("E19").Value = "Còn"
$ws.Range("F19").Value = "Tacos bò thơm ngon với salsa và guacamole."
# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Nachos Supreme"
$ws.Range("C20").Value = "Đồ Uống Nóng"
$ws.Range("D20").Value = 271201
$ws.Range("E20").Value = "Còn"
# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "Sua chua"
$ws.Range("C21").Value = "Đồ Uống Nóng"
$ws.Range("D21").Value = 10000
$ws.Range("E21").Value = "Còn"
$ws.Range("F21").Value = "do lanh"
# Row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "Them mon an"
$ws.Range("C22").Value = "Đồ Uống Nóng"
$ws.Range("D22").Value = 10000
$ws.Range("E22").Value = "Còn"
# Row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "Them mon an moi"
$ws.Range("C23").Value = "Đồ Uống Nóng"
$ws.Range("D23").Value = 23000
$ws.Range("E23").Value = "Còn"
$ws.Range("F23").Value = "ssdsddscss"
# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "sfdsfffdfsfee"
$ws.Range("C24").Value = "Đồ Uống Nóng"
$ws.Range("D24").Value = 10000
$ws.Range("E24").Value = "Còn"
# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "dfsdfdsfgg"
$ws.Range("C25").Value = "Đồ Uống Nóng"
$ws.Range("D25").Value = 444444
$ws.Range("E25").Value = "Còn"
# Row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "fdsfdsf32432"
$ws.Range("C26").Value = "Bia"
$ws.Range("D26").Value = 33243
$ws.Range("E26").Value = "Còn"
# Row 27
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "rtertre"
$ws.Range("C27").Value = "Đồ Uống Nóng"
$ws.Range("D27").Value = 4444
$ws.Range("E27").Value = "Còn"
# Row 28
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "fdfdsf"
$ws.Range("C28").Value = "Đồ Uống Nóng"
$ws.Range("D28").Value = 455
$ws.Range("E28").Value = "Còn"
# Row 29
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "fdfdsf"
$ws.Range("C29").Value = "Đồ Uống Nóng"
$ws.Range("D29").Value = 455
$ws.Range("E29").Value = "Còn"
# Row 30
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 34324
$ws.Range("C30").Value = "Bia"
$ws.Range("D30").Value = 444
$ws.Range("E30").Value = "Còn"
# Row 31
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 34324
$ws.Range("C31").Value = "Bia"
$ws.Range("D31").Value = 444
$ws.Range("E31").Value = "Còn"
# Row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = "efertre"
$ws.Range("C32").Value = "Đồ Uống Nóng"
$ws.Range("D32").Value = 3333
$ws.Range("E32").Value = "Còn"
# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "gdfsgdsf"
$ws.Range("C33").Value = "Đồ Uống Lạnh"
$ws.Range("D33").Value = 6456546
$ws.Range("E33").Value = "Còn"
# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 44
$ws.Range("C34").Value = "Đồ Uống Lạnh"
$ws.Range("D34").Value = 555
$ws.Range("E34").Value = "Còn"

# Apply formatting: center-align data cells, and thousands-separator number format on the price column (D)
$ws.Range("D3:D35").HorizontalAlignment = -4108
$ws.Range("D3:D35").NumberFormat = "#,##0"
$ws.Range("A3:C35").HorizontalAlignment = -4108
$ws.Range("E3:F35").HorizontalAlignment = -4108

# Update selection to mirror the final authored state
$ws.Range("A26:F35").Select()
